$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H95").Value = 66500
$ws.Range("J95").Value = 66500
$ws.Range("L95").Value = 66500
$ws.Range("N95").Value = -71992
$ws.Range("H132").Value = 3680777.5
$ws.Range("I132").Value = 4036594.8
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 12109784.4
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -12107254.4
$ws.Range("N132").Value = -17060
$ws.Range("H141").Value = 3226.3635
$ws.Range("I141").Value = 2972.0588
$ws.Range("J141").Value = 4091
$ws.Range("K141").Value = 8916.1764
$ws.Range("L141").Value = 12273
$ws.Range("M141").Value = -3736.1764
$ws.Range("N141").Value = -22633

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1315.3529
$ws.Range("I2").Value = 1257.4
$ws.Range("J2").Value = 1750
$ws.Range("K2").Value = 1257.4
$ws.Range("L2").Value = 1750
$ws.Range("M2").Value = -1144.4
$ws.Range("N2").Value = -1976
$ws.Range("H61").Value = 1688.5714
$ws.Range("I61").Value = 1688.5714
$ws.Range("K61").Value = 1688.5714
$ws.Range("M61").Value = -1476.5714
$ws.Range("H74").Value = 4578.077
$ws.Range("I74").Value = 3477.7778
$ws.Range("J74").Value = 7053.75
$ws.Range("K74").Value = 3477.7778
$ws.Range("L74").Value = 7053.75
$ws.Range("M74").Value = -2603.7778
$ws.Range("N74").Value = -8801.75
$ws.Range("H75").Value = 20890
$ws.Range("J75").Value = 20890
$ws.Range("L75").Value = 20890
$ws.Range("N75").Value = -22638
$ws.Range("H77").Value = 4578.077
$ws.Range("I77").Value = 3477.7778
$ws.Range("J77").Value = 7053.75
$ws.Range("K77").Value = 17388.889
$ws.Range("L77").Value = 35268.75
$ws.Range("M77").Value = -13020.889
$ws.Range("N77").Value = -44004.75
$ws.Range("H78").Value = 20890
$ws.Range("J78").Value = 20890
$ws.Range("L78").Value = 62670
$ws.Range("N78").Value = -71406
$ws.Range("H102").Value = 3652.4443
$ws.Range("I102").Value = 3546.1667
$ws.Range("J102").Value = 3865
$ws.Range("K102").Value = 3546.1667
$ws.Range("L102").Value = 3865
$ws.Range("M102").Value = -1924.1667
$ws.Range("N102").Value = -7109
$ws.Range("H116").Value = 1315.3529
$ws.Range("I116").Value = 1257.4
$ws.Range("J116").Value = 1750
$ws.Range("K116").Value = 1257.4
$ws.Range("L116").Value = 1750
$ws.Range("M116").Value = 1036.6
$ws.Range("N116").Value = -6338
$ws.Range("H132").Value = 1091.7
$ws.Range("I132").Value = 631.1429000000001
$ws.Range("J132").Value = 2166.3333
$ws.Range("K132").Value = 1893.4287
$ws.Range("L132").Value = 6498.999899999999
$ws.Range("M132").Value = 636.5712999999998
$ws.Range("N132").Value = -11558.9999
$ws.Range("H136").Value = 1688.5714
$ws.Range("I136").Value = 1688.5714
$ws.Range("K136").Value = 5065.7142
$ws.Range("M136").Value = -2515.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1315.3529
$ws.Range("I3").Value = 1257.4
$ws.Range("J3").Value = 1750
$ws.Range("K3").Value = 1257.4
$ws.Range("L3").Value = 1750
$ws.Range("M3").Value = -1143.4
$ws.Range("N3").Value = -1978
$ws.Range("H20").Value = 2512.9443
$ws.Range("J20").Value = 2401.3333
$ws.Range("L20").Value = 2401.3333
$ws.Range("N20").Value = -2895.3333
$ws.Range("H99").Value = 1591.5834
$ws.Range("I99").Value = 1429.9
$ws.Range("K99").Value = 1429.9
$ws.Range("M99").Value = 68.09999999999991
$ws.Range("H105").Value = 2510
$ws.Range("I105").Value = 2510
$ws.Range("K105").Value = 2510
$ws.Range("M105").Value = -763

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 51460.094
$ws.Range("I31").Value = 101872.3
$ws.Range("J31").Value = 5630.8184
$ws.Range("K31").Value = 101872.3
$ws.Range("L31").Value = 5630.8184
$ws.Range("M31").Value = -101577.3
$ws.Range("N31").Value = -6220.8184
$ws.Range("H34").Value = 51460.094
$ws.Range("I34").Value = 101872.3
$ws.Range("J34").Value = 5630.8184
$ws.Range("K34").Value = 101872.3
$ws.Range("L34").Value = 5630.8184
$ws.Range("M34").Value = -101670.3
$ws.Range("N34").Value = -6034.8184
$ws.Range("H98").Value = 219000
$ws.Range("J98").Value = 219000
$ws.Range("L98").Value = 219000
$ws.Range("N98").Value = -223492

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5263
$ws.Range("I70").Value = 5131.636
$ws.Range("J70").Value = 5443.625
$ws.Range("K70").Value = 5131.636
$ws.Range("L70").Value = 5443.625
$ws.Range("M70").Value = -4861.636
$ws.Range("N70").Value = -5983.625
$ws.Range("H73").Value = 5263
$ws.Range("I73").Value = 5131.636
$ws.Range("J73").Value = 5443.625
$ws.Range("K73").Value = 5131.636
$ws.Range("L73").Value = 5443.625
$ws.Range("M73").Value = -4195.636
$ws.Range("N73").Value = -7315.625
$ws.Range("H80").Value = 3666.6667
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 3800
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 3800
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -5796
$ws.Range("H83").Value = 3666.6667
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 3800
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 19000
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -28984
$ws.Range("H122").Value = 2546.3076
$ws.Range("I122").Value = 2484.8572
$ws.Range("J122").Value = 2618
$ws.Range("K122").Value = 7454.571599999999
$ws.Range("L122").Value = 7854
$ws.Range("M122").Value = -5004.571599999999
$ws.Range("N122").Value = -12754

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2228.7778
$ws.Range("I7").Value = 1291.5
$ws.Range("J7").Value = 4103.3335
$ws.Range("K7").Value = 1291.5
$ws.Range("L7").Value = 4103.3335
$ws.Range("M7").Value = -1179.5
$ws.Range("N7").Value = -4327.3335
$ws.Range("H46").Value = 576
$ws.Range("H122").Value = 4331.15
$ws.Range("I122").Value = 5134
$ws.Range("J122").Value = 3126.875
$ws.Range("K122").Value = 15402
$ws.Range("L122").Value = 9380.625
$ws.Range("M122").Value = -12952
$ws.Range("N122").Value = -14280.625
$ws.Range("H126").Value = 2228.7778
$ws.Range("I126").Value = 1291.5
$ws.Range("J126").Value = 4103.3335
$ws.Range("K126").Value = 3874.5
$ws.Range("L126").Value = 12310.0005
$ws.Range("M126").Value = -1404.5
$ws.Range("N126").Value = -17250.0005
$ws.Range("H136").Value = 4183.3228
$ws.Range("I136").Value = 5083.8486
$ws.Range("J136").Value = 3158.5862
$ws.Range("K136").Value = 15251.5458
$ws.Range("L136").Value = 9475.758600000001
$ws.Range("M136").Value = -12701.5458
$ws.Range("N136").Value = -14575.7586

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 37514.5
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 37514.5
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 37514.5
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -38148.5
$ws.Range("H38").Value = 70000
$ws.Range("J38").Value = 70000
$ws.Range("L38").Value = 70000
$ws.Range("N38").Value = -70946
$ws.Range("H49").Value = 37500
$ws.Range("I49").Value = 5000
$ws.Range("J49").Value = 70000
$ws.Range("K49").Value = 5000
$ws.Range("L49").Value = 70000
$ws.Range("M49").Value = -4770
$ws.Range("N49").Value = -70460
$ws.Range("H125").Value = 40276.07
$ws.Range("J125").Value = 40276.07
$ws.Range("L125").Value = 40276.07
$ws.Range("N125").Value = -50116.07
$ws.Range("H131").Value = 38282.367
$ws.Range("J131").Value = 38282.367
$ws.Range("L131").Value = 38282.367
$ws.Range("N131").Value = -48362.367
